# === casos_teste.xlsx edit script ===
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (Planilha1 -> casos_teste)
$ws.Name = 'casos_teste'

# --- Fix up cell formatting (fill/border/numberformat) before writing new values ---
# Row 3 becomes a copy of the green "values" row style (style index 2)
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A3:E3").PasteSpecial(-4122) | Out-Null

# Rows 7 & 8 column E needs the pink row style (style index 3) instead of the gray "Ruim" style
$ws.Range("A7").Copy() | Out-Null
$ws.Range("E7:E8").PasteSpecial(-4122) | Out-Null

# Rows 9-13 (new rows) need the plain pink row style (style index 3) across A:C and E
$ws.Range("A7:C7").Copy() | Out-Null
$ws.Range("A9:C13").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Copy() | Out-Null
$ws.Range("E9:E13").PasteSpecial(-4122) | Out-Null
$ws.Range("C7").Copy() | Out-Null
$ws.Range("D9:D13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Write the final cell values (rows 3-13; rows 1-2 are unchanged) ---
$ws.Range("A3").Value = 'juros_composto.py'
$ws.Range("B3").Value = 'calcular_juros_compostos'
$ws.Range("C3").Value = 'Enviando o juros como número decimal'
$ws.Range("D3").Value = '1000, 15.5, 1'
$ws.Range("E3").Value = '(155, 1155)'

$ws.Range("A4").Value = 'juros_composto.py'
$ws.Range("B4").Value = 'calcular_juros_compostos'
$ws.Range("C4").Value = 'Enviando o capital como número negativo'
$ws.Range("D4").Value = '-1000, 40, 2'
$ws.Range("E4").Value = 'ValueError("O capital investido não pode ser negativo.")'

$ws.Range("A5").Value = 'juros_composto.py'
$ws.Range("B5").Value = 'calcular_juros_compostos'
$ws.Range("C5").Value = 'Enviando o juros como número negativo'
$ws.Range("D5").Value = '1000, -40, 2'
$ws.Range("E5").Value = 'ValueError("O capital investido não pode ser negativo.")'

$ws.Range("A6").Value = 'juros_composto.py'
$ws.Range("B6").Value = 'calcular_juros_compostos'
$ws.Range("C6").Value = 'Enviando o tempo como número negativo'
$ws.Range("D6").Value = '1000, 40, -2'
$ws.Range("E6").Value = 'ValueError("O capital investido não pode ser negativo.")'

$ws.Range("A7").Value = 'juros_composto.py'
$ws.Range("B7").Value = 'calcular_juros_compostos'
$ws.Range("C7").Value = 'Enviando algo que não seja um número no lugar do capital '
$ws.Range("D7").Value = 'ola, 40, 2'
$ws.Range("E7").Value = 'TypeError("O capital investido deve ser um número (int ou float).")'

$ws.Range("A8").Value = 'juros_composto.py'
$ws.Range("B8").Value = 'calcular_juros_compostos'
$ws.Range("C8").Value = 'Enviando algo que não seja um número no lugar dos juros'
$ws.Range("D8").Value = '1000, "ola", 2'
$ws.Range("E8").Value = 'TypeError("A taxa de juros deve ser um número (int ou float).")'

$ws.Range("A9").Value = 'juros_composto.py'
$ws.Range("B9").Value = 'calcular_juros_compostos'
$ws.Range("C9").Value = 'Enviando algo que não seja um número no lugar do tempo '
$ws.Range("D9").Value = '1000, 40, "ola"'
$ws.Range("E9").Value = 'TypeError("O tempo deve ser um número (int ou float).")'

$ws.Range("A10").Value = 'juros_composto.py'
$ws.Range("B10").Value = 'calcular_juros_compostos'
$ws.Range("C10").Value = 'Enviando menos que 3 valores'
$ws.Range("D10").Value = '3500, 40'
$ws.Range("E10").Value = 'ValueError("Não é permitido enviar menos que 3 valores")'

$ws.Range("A11").Value = 'juros_composto.py'
$ws.Range("B11").Value = 'calcular_juros_compostos'
$ws.Range("C11").Value = 'Enviando o número 0 no lugar do capital'
$ws.Range("D11").Value = '0, 40, 2'
$ws.Range("E11").Value = 'ValueError("O capital deve ser um número maior que 0")'

$ws.Range("A12").Value = 'juros_composto.py'
$ws.Range("B12").Value = 'calcular_juros_compostos'
$ws.Range("C12").Value = 'Enviando o número 0 no lugar do juros'
$ws.Range("D12").Value = '1000, 0, 2'
$ws.Range("E12").Value = 'ValueError("O juros deve ser um número maior que 0")'

$ws.Range("A13").Value = 'juros_composto.py'
$ws.Range("B13").Value = 'calcular_juros_compostos'
$ws.Range("C13").Value = 'Enviando o número 0 no lugar do tempo'
$ws.Range("D13").Value = '1000, 40, 0'
$ws.Range("E13").Value = 'ValueError("O tempo deve ser um número maior que 0")'

# --- Column widths (approximate: engine quantizes to 1/6-character steps) ---
$ws.Columns.Item(3).ColumnWidth = 56.666666666666664
$ws.Columns.Item(5).ColumnWidth = 63.0

# --- Row heights for the newly-added rows (match the 27.95pt custom height used throughout) ---
$ws.Rows.Item(9).RowHeight = 27.95
$ws.Rows.Item(10).RowHeight = 27.95
$ws.Rows.Item(11).RowHeight = 27.95
$ws.Rows.Item(12).RowHeight = 27.95
$ws.Rows.Item(13).RowHeight = 27.95

# --- View state: zoom + selection ---
$ws.Range("C8").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
